# Apply the "Add files via upload" edit:
#  - Mark three lab results as "pass" (text) instead of blank:
#      Sheet "БИВТ-22-17" (sheet1): E18  -> "pass"  (3rd lab, Крамаренко А.В.)
#      Sheet "БИВТ-22-18" (sheet2): H3   -> "pass"  (6th lab, Гафаров Э.Х.)
#      Sheet "БИВТ-22-18" (sheet2): G23  -> "pass"  (5th lab, Сентяков Г.Д.)
#  - Update the per-column summary-row formulas from COUNT(...) to
#    COUNTA(...) on every sheet, so the newly-added text "pass" marks
#    are reflected in the totals.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("БИВТ-22-17")
$sheet2 = $wb.Worksheets.Item("БИВТ-22-18")
$sheet3 = $wb.Worksheets.Item("БИВТ-22-20")

# --- Set the newly-filled-in "pass" marks ---------------------------------
$sheet1.Range("E18").Value = "pass"
$sheet2.Range("H3").Value = "pass"
$sheet2.Range("G23").Value = "pass"

# --- Switch the summary row from COUNT to COUNTA --------------------------
# Each column keeps its own relative-reference formula (no $ anchors),
# matching the original COUNT(...) formulas that are being replaced.

# Sheet1 ("БИВТ-22-17"): summary row 33, columns B:L (data rows 4:32)
$cols = @("B","C","D","E","F","G","H","I","J","K","L")
foreach ($col in $cols) {
    $sheet1.Range($col + "33").Formula = "=COUNTA(" + $col + "4:" + $col + "32)"
}

# Sheet2 ("БИВТ-22-18"): summary row 27, columns B:L (data rows 2:26)
foreach ($col in $cols) {
    $sheet2.Range($col + "27").Formula = "=COUNTA(" + $col + "2:" + $col + "26)"
}

# Sheet3 ("БИВТ-22-20"): summary row 31, columns B:L (data rows 2:30)
foreach ($col in $cols) {
    $sheet3.Range($col + "31").Formula = "=COUNTA(" + $col + "2:" + $col + "30)"
}

$excel.Calculate()

# --- Leave the workbook's view state the way the author saved it ----------
[void]$sheet1.Range("E19").Select()
[void]$sheet2.Range("G24").Select()
[void]$sheet3.Select()
[void]$sheet3.Range("B31:L31").Select()
